$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 83
$ws.Range("B1").Value = 83
$ws.Range("C1").Value = 169

$ws.Range("A2").Value = 83
$ws.Range("B2").Value = 192.2000000000007
$ws.Range("C2").Value = 83
